$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Sending: FAPs, Target: ECs)
$ws.Range("M2").Value = 1.028010333333333
$ws.Range("N2").Value = 3.084031
$ws.Range("O2").Value = 0.04165745457248912
$ws.Range("P2").Value = 0.04165745457248914
$ws.Range("Q2").Value = 2.147571497582111
$ws.Range("R2").Value = 19.328143478239
$ws.Range("S2").Value = 0.0299792200841283
$ws.Range("T2").Value = 0.02997922008412832

# Row 3 (Sending: FAPs, Target: FAPs)
$ws.Range("O3").Value = 0.9361395479363341
$ws.Range("P3").Value = 0.9361395479363344
$ws.Range("S3").Value = 0.6737025539619477
$ws.Range("T3").Value = 0.6737025539619479

# Row 4 (Sending: FAPs, Target: MuSCs)
$ws.Range("O4").Value = 0.02220299749117665
$ws.Range("P4").Value = 0.02220299749117666
$ws.Range("S4").Value = 0.01597861787635289
$ws.Range("T4").Value = 0.0159786178763529

# Row 5 (Sending: MuSCs, Target: ECs)
$ws.Range("G5").Value = 0.8137799999999999
$ws.Range("I5").Value = 0.280339608077571
$ws.Range("J5").Value = 0.280339608077571
$ws.Range("M5").Value = 1.028010333333333
$ws.Range("N5").Value = 3.084031
$ws.Range("O5").Value = 0.04165745457248912
$ws.Range("P5").Value = 0.04165745457248914
$ws.Range("Q5").Value = 0.8365742490599999
$ws.Range("R5").Value = 7.52916824154
$ws.Range("S5").Value = 0.01167823448836082
$ws.Range("T5").Value = 0.01167823448836082

# Row 6 (Sending: MuSCs, Target: FAPs)
$ws.Range("G6").Value = 0.8137799999999999
$ws.Range("I6").Value = 0.280339608077571
$ws.Range("J6").Value = 0.280339608077571
$ws.Range("O6").Value = 0.9361395479363341
$ws.Range("P6").Value = 0.9361395479363344
$ws.Range("S6").Value = 0.2624369939743864
$ws.Range("T6").Value = 0.2624369939743865

# Row 7 (Sending: MuSCs, Target: MuSCs)
$ws.Range("G7").Value = 0.8137799999999999
$ws.Range("I7").Value = 0.280339608077571
$ws.Range("J7").Value = 0.280339608077571
$ws.Range("O7").Value = 0.02220299749117665
$ws.Range("P7").Value = 0.02220299749117666
$ws.Range("Q7").Value = 0.44588552382
$ws.Range("S7").Value = 0.006224379614823755
$ws.Range("T7").Value = 0.006224379614823757
